# Mark the newly-added "隐藏控件" (hideControlElement) and "显示控件"
# (showControlElement) function rows as completed on the "事件功能" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("事件功能")

# B1 already carries the "已完成" (completed) formatting - reuse it so the
# newly finished rows (11 & 12) get the same green fill/border style.
$template = $ws.Range("B1")
$template.Copy()

$ws.Range("B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B11").Value = "已完成"
$ws.Range("B12").Value = "已完成"

$excel.CutCopyMode = 0

$ws.Activate()
$ws.Range("B15").Select()
